# Daily attendance processing - swap the order of the first two
# "Recorded By" entries (column G) for every data row on the active sheet.
# Entries are a comma-separated list (e.g. "System, user@example.com");
# rows with a single entry are left untouched, and any extra entries
# beyond the first two (e.g. a trailing "system" tag) keep their position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Value2

    if ($text -and $text.Contains(",")) {
        $parts = $text.Split(",")
        if ($parts.Length -ge 2) {
            $trimmed = @()
            foreach ($p in $parts) {
                $trimmed += $p.Trim()
            }

            $first = $trimmed[0]
            $second = $trimmed[1]
            $trimmed[0] = $second
            $trimmed[1] = $first

            $newText = [string]::Join(", ", $trimmed)
            $cell.Value2 = $newText
        }
    }
}
